$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update values
$ws.Range("B3").Value = 0.01474759464883743
$ws.Range("C3").Value = 0.01511028007094146
$ws.Range("D3").Value = 257891443789675.2

# Row 4 - rename to DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01552429973541219
$ws.Range("C4").Value = 0.01658962487133521
$ws.Range("D4").Value = 74684523125812.67

# Row 5 - rename to MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 49106809746385.05
$ws.Range("C5").Value = 22001254883931.89
$ws.Range("D5").Value = 246334403328739.3
